# Updated cryptos list on Tue Feb 28 18:00:34 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (D) and "Volume(1h)" (E) columns for every coin row
# (rows 2-51) on the active sheet with the latest scraped values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> @{ D = "<new price text>" (optional); E = "<new volume text>" }
$updates = [ordered]@{
    2  = @{ D = "23.526.45";    E = "  +0.58%  " }
    3  = @{ D = "1.642.24";     E = "  +0.29%  " }
    4  = @{ D = "0.9989";       E = "  -0.49%  " }
    5  = @{ D = "0.9991";       E = "  -0.31%  " }
    6  = @{ D = "303.74";       E = "  -0.06%  " }
    7  = @{ D = "0.3795";       E = "  +0.33%  " }
    8  = @{ D = "51.97";        E = "  -0.61%  " }
    9  = @{ D = "0.3629";       E = "  -0.45%  " }
    10 = @{             E = "  +1.20%  " }
    11 = @{             E = "  -0.50%  " }
    12 = @{ D = "0.9988";       E = "  -0.55%  " }
    13 = @{ D = "22.63";        E = "  -0.82%  " }
    14 = @{ D = "6.482";        E = "  -2.23%  " }
    15 = @{ D = "7.388";        E = "  +1.60%  " }
    16 = @{ D = "0.00001242";   E = "  -0.57%  " }
    17 = @{ D = "1.636.54";     E = "  -0.15%  " }
    18 = @{ D = "95.32";        E = "  +1.39%  " }
    19 = @{ D = "0.06945";      E = "  +0.06%  " }
    20 = @{ D = "17.60";        E = "  -2.85%  " }
    21 = @{ D = "6.573";        E = "  +0.61%  " }
    22 = @{ D = "0.9992";       E = "  -0.30%  " }
    23 = @{             E = "  -2.36%  " }
    24 = @{ D = "23.521.41";    E = "  +0.47%  " }
    25 = @{ D = "2.523";        E = "  +2.93%  " }
    26 = @{ D = "3.080";        E = "  -5.24%  " }
    27 = @{ D = "21.24";        E = "  +0.18%  " }
    28 = @{ D = "152.51";       E = "  +2.27%  " }
    29 = @{ D = "5.261";        E = "  -0.82%  " }
    30 = @{             E = "  -1.72%  " }
    31 = @{ D = "1.817.78";     E = "  -0.21%  " }
    32 = @{ D = "1.103";        E = "  +14.56%  " }
    33 = @{ D = "6.634";        E = "  -3.39%  " }
    34 = @{             E = "  -7.09%  " }
    35 = @{ D = "11.48";        E = "  +4.57%  " }
    36 = @{ D = "0.02772";      E = "  -3.27%  " }
    37 = @{ D = "0.2511";       E = "  -1.79%  " }
    38 = @{             E = "  -1.29%  " }
    39 = @{ D = "6.043";        E = "  -3.27%  " }
    40 = @{ D = "0.07071";      E = "  -2.58%  " }
    41 = @{ D = "0.7081";       E = "  -0.32%  " }
    42 = @{ D = "1.355";        E = "  -1.41%  " }
    43 = @{ D = "12.33";        E = "  -1.58%  " }
    44 = @{ D = "15.62";        E = "  -4.72%  " }
    45 = @{             E = "  +0.36%  " }
    46 = @{ D = "0.9984";       E = "  -0.30%  " }
    47 = @{ D = "2.294";        E = "  -2.42%  " }
    48 = @{             E = "  -0.67%  " }
    49 = @{             E = "  -0.01%  " }
    50 = @{ D = "129.28";       E = "  +1.43%  " }
    51 = @{ D = "1.200";        E = "  -1.63%  " }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]

    if ($vals.ContainsKey("D")) {
        $priceCell = $ws.Range("D$row")
        $priceText = $vals.D

        # Several prices ("0.9989", "22.63", "1.200", ...) are plain decimal
        # literals that Excel's COM layer would otherwise auto-convert to a
        # Number (and, worse, normalize away trailing zeros, e.g. "1.200"
        # -> 1.2). Force the cell to Text first so the exact digit string -
        # including any trailing zeros - round-trips, then restore the
        # default "Normal" style so no stray formatting is left behind.
        $looksNumeric = $priceText -match '^[+-]?\d+(\.\d+)?$'
        if ($looksNumeric) {
            $priceCell.NumberFormat = "@"
            $priceCell.Value = $priceText
            $priceCell.Style = "Normal"
        } else {
            $priceCell.Value = $priceText
        }
    }

    $ws.Range("E$row").Value = $vals.E
}
